# Auto-generated: apply scheduled market-price data refresh to Leve profit tables
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
  # ALC row 9 (hunk 0)
  $ws.Range("H9").Value = 323.58334
  $ws.Range("I9").Value = 269
  $ws.Range("J9").Value = 362.57144
  $ws.Range("K9").Value = 269
  $ws.Range("L9").Value = 362.57144
  $ws.Range("M9").Value = -100
  $ws.Range("N9").Value = -700.5714399999999
  # ALC row 53 (hunk 1)
  $ws.Range("H53").Value = 277.91666
  $ws.Range("I53").Value = 107.25
  $ws.Range("J53").Value = 363.25
  $ws.Range("K53").Value = 107.25
  $ws.Range("L53").Value = 363.25
  $ws.Range("M53").Value = 529.75
  $ws.Range("N53").Value = -1637.25
  # ALC row 136 (hunk 2)
  $ws.Range("H136").Value = 39750
  $ws.Range("J136").Value = 39750
  $ws.Range("L136").Value = 39750
  $ws.Range("N136").Value = -49950
  # ALC row 138 (hunk 3)
  $ws.Range("H138").Value = 2815.889
  $ws.Range("I138").Value = 1248.8572
  $ws.Range("J138").Value = 4187.0415
  $ws.Range("K138").Value = 3746.5716
  $ws.Range("L138").Value = 12561.1245
  $ws.Range("M138").Value = 1393.4284
  $ws.Range("N138").Value = -22841.1245

$ws = $wb.Worksheets.Item("ARM")
  # ARM row 32 (hunk 4)
  $ws.Range("H32").Value = 4986.2104
  $ws.Range("I32").Value = 5728.0215
  $ws.Range("K32").Value = 5728.0215
  $ws.Range("M32").Value = -5441.0215
  # ARM row 45 (hunk 5)
  $ws.Range("H45").Value = 1913.3864
  $ws.Range("I45").Value = 1639.7428
  $ws.Range("J45").Value = 2977.5557
  $ws.Range("K45").Value = 1639.7428
  $ws.Range("L45").Value = 2977.5557
  $ws.Range("M45").Value = -1262.7428
  $ws.Range("N45").Value = -3731.5557
  # ARM row 102 (hunk 6)
  $ws.Range("H102").Value = 4164.727
  $ws.Range("I102").Value = 2145.7778
  $ws.Range("J102").Value = 13250
  $ws.Range("K102").Value = 2145.7778
  $ws.Range("L102").Value = 13250
  $ws.Range("M102").Value = -523.7777999999998
  $ws.Range("N102").Value = -16494

$ws = $wb.Worksheets.Item("BSM")
  # BSM row 20 (hunk 7)
  $ws.Range("H20").Value = 1322.0625
  $ws.Range("I20").Value = 1610.5
  $ws.Range("K20").Value = 1610.5
  $ws.Range("M20").Value = -1363.5
  # BSM row 132 (hunk 8)
  $ws.Range("H132").Value = 55565.633
  $ws.Range("J132").Value = 55565.633
  $ws.Range("L132").Value = 55565.633
  $ws.Range("N132").Value = -65685.633

$ws = $wb.Worksheets.Item("CRP")
  # CRP row 31 (hunk 9)
  $ws.Range("H31").Value = 1938.9688
  $ws.Range("I31").Value = 1497.5238
  $ws.Range("J31").Value = 2781.7273
  $ws.Range("K31").Value = 1497.5238
  $ws.Range("L31").Value = 2781.7273
  $ws.Range("M31").Value = -1202.5238
  $ws.Range("N31").Value = -3371.7273
  # CRP row 34 (hunk 10)
  $ws.Range("H34").Value = 1938.9688
  $ws.Range("I34").Value = 1497.5238
  $ws.Range("J34").Value = 2781.7273
  $ws.Range("K34").Value = 1497.5238
  $ws.Range("L34").Value = 2781.7273
  $ws.Range("M34").Value = -1295.5238
  $ws.Range("N34").Value = -3185.7273
  # CRP row 86 (hunk 11)
  $ws.Range("H86").Value = 31731.312
  $ws.Range("I86").Value = 43847.625
  $ws.Range("J86").Value = 19615
  $ws.Range("K86").Value = 43847.625
  $ws.Range("L86").Value = 19615
  $ws.Range("M86").Value = -42724.625
  $ws.Range("N86").Value = -21861
  # CRP row 89 (hunk 12)
  $ws.Range("H89").Value = 31731.312
  $ws.Range("I89").Value = 43847.625
  $ws.Range("J89").Value = 19615
  $ws.Range("K89").Value = 219238.125
  $ws.Range("L89").Value = 98075
  $ws.Range("M89").Value = -213622.125
  $ws.Range("N89").Value = -109307
  # CRP row 94 (hunk 13)
  $ws.Range("H94").Value = 1531.1
  $ws.Range("J94").Value = 1847.8334
  $ws.Range("L94").Value = 1847.8334
  $ws.Range("N94").Value = -2749.8334

$ws = $wb.Worksheets.Item("CUL")
  # CUL row 51 (hunk 14)
  $ws.Range("H51").Value = 169.66667
  $ws.Range("I51").Value = 169.66667
  $ws.Range("K51").Value = 509.00001
  $ws.Range("M51").Value = -49.00001000000003
  # CUL row 63 (hunk 15)
  $ws.Range("H63").Value = 2750
  $ws.Range("I63").Value = 2750
  $ws.Range("J63").Value = 0
  $ws.Range("K63").Value = 8250
  $ws.Range("L63").Value = 0
  $ws.Range("M63").Value = -7501
  $ws.Range("N63").ClearContents()
  # CUL row 66 (hunk 16)
  $ws.Range("H66").Value = 2750
  $ws.Range("I66").Value = 2750
  $ws.Range("J66").Value = 0
  $ws.Range("K66").Value = 24750
  $ws.Range("L66").Value = 0
  $ws.Range("M66").Value = -21006
  $ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
  # GSM row 52 (hunk 17)
  $ws.Range("H52").Value = 25822
  $ws.Range("J52").Value = 29983
  $ws.Range("L52").Value = 29983
  $ws.Range("N52").Value = -30501
  # GSM row 102 (hunk 18)
  $ws.Range("H102").Value = 3583.2593
  $ws.Range("I102").Value = 2581.7273
  $ws.Range("J102").Value = 7990
  $ws.Range("K102").Value = 2581.7273
  $ws.Range("L102").Value = 7990
  $ws.Range("M102").Value = -959.7273
  $ws.Range("N102").Value = -11234
  # GSM row 122 (hunk 19)
  $ws.Range("H122").Value = 2721.1304
  $ws.Range("I122").Value = 2095
  $ws.Range("J122").Value = 3295.0833
  $ws.Range("K122").Value = 6285
  $ws.Range("L122").Value = 9885.249899999999
  $ws.Range("M122").Value = -3835
  $ws.Range("N122").Value = -14785.2499
  # GSM row 126 (hunk 20)
  $ws.Range("H126").Value = 2001
  $ws.Range("I126").Value = 2001
  $ws.Range("K126").Value = 6003
  $ws.Range("M126").Value = -3533
  # GSM row 132 (hunk 21)
  $ws.Range("H132").Value = 6416.878
  $ws.Range("I132").Value = 5575.8486
  $ws.Range("J132").Value = 9886.125
  $ws.Range("K132").Value = 16727.5458
  $ws.Range("L132").Value = 29658.375
  $ws.Range("M132").Value = -14197.5458
  $ws.Range("N132").Value = -34718.375

$ws = $wb.Worksheets.Item("LTW")
  # LTW row 7 (hunk 22)
  $ws.Range("H7").Value = 2584.6365
  $ws.Range("I7").Value = 2584.6365
  $ws.Range("K7").Value = 2584.6365
  $ws.Range("M7").Value = -2472.6365
  # LTW row 16 (hunk 23)
  $ws.Range("H16").Value = 22727566
  $ws.Range("I16").Value = 35714508
  $ws.Range("J16").Value = 419.5
  $ws.Range("K16").Value = 35714508
  $ws.Range("L16").Value = 419.5
  $ws.Range("M16").Value = -35714338
  $ws.Range("N16").Value = -759.5
  # LTW row 43 (hunk 24)
  $ws.Range("H43").Value = 32852.438
  $ws.Range("I43").Value = 16115.25
  $ws.Range("J43").Value = 38431.5
  $ws.Range("K43").Value = 16115.25
  $ws.Range("L43").Value = 38431.5
  $ws.Range("M43").Value = -15922.25
  $ws.Range("N43").Value = -38817.5
  # LTW row 100 (hunk 25)
  $ws.Range("H100").Value = 6197.5
  $ws.Range("I100").Value = 5162.5
  $ws.Range("J100").Value = 7750
  $ws.Range("K100").Value = 5162.5
  $ws.Range("L100").Value = 7750
  $ws.Range("M100").Value = -4621.5
  $ws.Range("N100").Value = -8832
  # LTW row 126 (hunk 26)
  $ws.Range("H126").Value = 2584.6365
  $ws.Range("I126").Value = 2584.6365
  $ws.Range("K126").Value = 7753.9095
  $ws.Range("M126").Value = -5283.9095

$ws = $wb.Worksheets.Item("WVR")
  # WVR row 122 (hunk 27)
  $ws.Range("H122").Value = 3705.3
  $ws.Range("I122").Value = 3439
  $ws.Range("K122").Value = 10317
  $ws.Range("M122").Value = -7867
  # WVR row 136 (hunk 28)
  $ws.Range("H136").Value = 3099.3333
  $ws.Range("I136").Value = 3649
  $ws.Range("K136").Value = 10947
  $ws.Range("M136").Value = -8397
